# Add a new weekly price record.
# A new row is inserted at row 36 (pushing the existing rows 36-114 down to
# 37-115) and populated with a fresh data point; all other rows keep their
# original values, just shifted down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:36").Insert()

$ws.Range("A36").Value = 4
$ws.Range("B36").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C36").Value = "Los Lagos"
$ws.Range("D36").Value = "2021-10-29"
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100101
$ws.Range("H36").Value = "Berries"
$ws.Range("I36").Value = 100112025
$ws.Range("J36").Value = "Frutilla"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 800
$ws.Range("N36").Value = 9500
$ws.Range("O36").Value = 10000
$ws.Range("P36").Value = 9750
$ws.Range("Q36").Value = "`$/bandeja 7 kilos"
$ws.Range("R36").Value = "Provincia de Melipilla"
$ws.Range("S36").Value = 1393
$ws.Range("T36").Value = 7
